$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0.22914666016441
$ws.Range("C2").Value = 0.22914666016441
$ws.Range("D2").Value = 1.15476288483359
$ws.Range("E2").Value = 0.00784726815629703
$ws.Range("F2").Value = 0.3225

# Row 3 updates
$ws.Range("B3").Value = 28.9716727333378
$ws.Range("C3").Value = 0.198436114611902
$ws.Range("E3").Value = 0.992152731843703

# Row 4 updates
$ws.Range("B4").Value = 29.2008193935022
